$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 4")

# Row 9: set Activity (column G) to "Prep."
$ws.Range("G9").Value = "Prep."

# Row 10: fill in Date, Start, Stop, Activity, Comments
$ws.Range("B10").Value = 43517
$ws.Range("C10").Value = 0.67708333333333337
$ws.Range("D10").Value = 0.73611111111111116
$ws.Range("G10").Value = "Prep."
$ws.Range("H10").Value = "Completing homework on UI/UX"

# Row 11: fill in Date, Start, Stop, Activity
$ws.Range("B11").Value = 43518
$ws.Range("C11").Value = 0.33333333333333331
$ws.Range("D11").Value = 0.41666666666666669
$ws.Range("G11").Value = "Class"

# Row 12: fill in Start time
$ws.Range("C12").Value = 0.70138888888888884

# Update the selected / active cell on the sheet
$ws.Activate()
$ws.Range("C13").Select()
